# Apply updated dSF (column F) values as part of a data repull / mean calculation update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -8
    3  = -3
    6  = -3
    8  = -2
    10 = 7
    11 = -1
    12 = 0
    18 = 12
    20 = 1
    22 = 1
    23 = 3
    24 = -1
    25 = -1
    28 = 6
    30 = -4
    32 = -1
    33 = -1
    34 = -2
    36 = 9
    37 = 2
    39 = 2
    44 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
